# Auto-generated script to update pl_mw.xlsx values for Case_5_36 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3429493964269454
$ws.Range("C2").Value = 0.06965172002794873
$ws.Range("D2").Value = 0.02613535110029375
$ws.Range("E2").Value = 0.4172994359769717
$ws.Range("F2").Value = 0.6289244698071883
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.4725067101974112
$ws.Range("K2").Value = 0.3790598273342312
$ws.Range("O2").Value = 2.123885881344393
# Row 3
$ws.Range("B3").Value = 0.3009604902053695
$ws.Range("C3").Value = 0.061536622470868
$ws.Range("D3").Value = 0.02409310344494742
$ws.Range("E3").Value = 0.3641208024745595
$ws.Range("F3").Value = 0.6278367374468701
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.4777763309126506
$ws.Range("K3").Value = 0.3309631340923715
$ws.Range("O3").Value = 2.133867296146363
# Row 4
$ws.Range("B4").Value = 0.2751380727278843
$ws.Range("C4").Value = 0.05652751646390186
$ws.Range("D4").Value = 0.02282952962679019
$ws.Range("E4").Value = 0.3315522843217025
$ws.Range("F4").Value = 0.6276358551096948
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.4813320661461979
$ws.Range("K4").Value = 0.301355535434908
$ws.Range("O4").Value = 2.141531551795708
# Row 5
$ws.Range("B5").Value = 0.2646054256215109
$ws.Range("C5").Value = 0.05447967130373854
$ws.Range("D5").Value = 0.02231222978210212
$ws.Range("E5").Value = 0.3182997360578668
$ws.Range("F5").Value = 0.6276712106583346
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.4828614232459323
$ws.Range("K5").Value = 0.2892715992567787
$ws.Range("O5").Value = 2.145040146666858
# Row 6
$ws.Range("B6").Value = 0.2628559110281117
$ws.Range("C6").Value = 0.05413923129175657
$ws.Range("D6").Value = 0.02222618969180701
$ws.Range("E6").Value = 0.3161002854831452
$ws.Range("F6").Value = 0.6276841548369632
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.4831202215763746
$ws.Range("K6").Value = 0.287263962931803
$ws.Range("O6").Value = 2.145645996966252
# Row 7
$ws.Range("B7").Value = 0.2749960647586533
$ws.Range("C7").Value = 0.05649992514496205
$ws.Range("D7").Value = 0.02282256274050098
$ws.Range("E7").Value = 0.3313734793942444
$ws.Range("F7").Value = 0.62763585760392
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.4813523663824455
$ws.Range("K7").Value = 0.3011926418602684
$ws.Range("O7").Value = 2.141577310738967
# Row 8
$ws.Range("B8").Value = 0.3284804428742234
$ws.Range("C8").Value = 0.06685913912150454
$ws.Range("D8").Value = 0.02543320350022071
$ws.Range("E8").Value = 0.3989451711657068
$ws.Range("F8").Value = 0.6284523662408716
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.4742571391747532
$ws.Range("K8").Value = 0.3624921648032284
$ws.Range("O8").Value = 2.127008269994576
# Row 9
$ws.Range("B9").Value = 0.4330196250668905
$ws.Range("C9").Value = 0.08696300939867285
$ws.Range("D9").Value = 0.03047489210479881
$ws.Range("E9").Value = 0.5321920092805641
$ws.Range("F9").Value = 0.6337695595010544
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.4628902603821636
$ws.Range("K9").Value = 0.482080110783528
$ws.Range("O9").Value = 2.110660902985103
# Row 10
$ws.Range("B10").Value = 0.5095990029419966
$ws.Range("C10").Value = 0.1016050693617387
$ws.Range("D10").Value = 0.03413008150781138
$ws.Range("E10").Value = 0.6306550029007525
$ws.Range("F10").Value = 0.6399578769489622
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.4561006330277806
$ws.Range("K10").Value = 0.5695492469239412
$ws.Range("O10").Value = 2.106157102289842
# Row 11
$ws.Range("B11").Value = 0.5443850744516681
$ws.Range("C11").Value = 0.1082383929685591
$ws.Range("D11").Value = 0.03578199605599508
$ws.Range("E11").Value = 0.6755970522736021
$ws.Range("F11").Value = 0.643272201252536
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.4533528994753127
$ws.Range("K11").Value = 0.6092537054502145
$ws.Range("O11").Value = 2.105750366634226
# Row 12
$ws.Range("B12").Value = 0.5575500160009597
$ws.Range("C12").Value = 0.1107462947899194
$ws.Range("D12").Value = 0.03640594065366543
$ws.Range("E12").Value = 0.6926390146231967
$ws.Range("F12").Value = 0.6445993015316844
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.4523615882685768
$ws.Range("K12").Value = 0.6242760389982038
$ws.Range("O12").Value = 2.105833419773376
# Row 13
$ws.Range("B13").Value = 0.5547150655618225
$ws.Range("C13").Value = 0.1102063516832743
$ws.Range("D13").Value = 0.03627163478749651
$ws.Range("E13").Value = 0.6889676557750732
$ws.Range("F13").Value = 0.6443102788797646
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.4525728940902773
$ws.Range("K13").Value = 0.621041292549279
$ws.Range("O13").Value = 2.105804973691676
# Row 14
$ws.Range("B14").Value = 0.5454683214053091
$ws.Range("C14").Value = 0.1084447998897531
$ws.Range("D14").Value = 0.03583336063724119
$ws.Range("E14").Value = 0.6769986262960117
$ws.Range("F14").Value = 0.6433799374226794
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.4532703567462519
$ws.Range("K14").Value = 0.6104898632670483
$ws.Range("O14").Value = 2.105752442689806
# Row 15
$ws.Range("B15").Value = 0.5398033927484391
$ws.Range("C15").Value = 0.1073652769510147
$ws.Range("D15").Value = 0.03556469558551356
$ws.Range("E15").Value = 0.6696703463099425
$ws.Range("F15").Value = 0.6428194652369967
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.4537039847019884
$ws.Range("K15").Value = 0.6040251138252017
$ws.Range("O15").Value = 2.105751168284797
# Row 16
$ws.Range("B16").Value = 0.5073245918154896
$ws.Range("C16").Value = 0.1011710099069205
$ws.Range("D16").Value = 0.03402190330881183
$ws.Range("E16").Value = 0.6277211131709919
$ws.Range("F16").Value = 0.6397513451047416
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.4562870786100106
$ws.Range("K16").Value = 0.5669526948046268
$ws.Range("O16").Value = 2.106216806973435
# Row 17
$ws.Range("B17").Value = 0.4873866239897495
$ws.Range("C17").Value = 0.097363967834184
$ws.Range("D17").Value = 0.03307264324293868
$ws.Range("E17").Value = 0.602026435657379
$ws.Range("F17").Value = 0.6379971943104934
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.4579591603172055
$ws.Range("K17").Value = 0.5441876184393379
$ws.Range("O17").Value = 2.106923707552085
# Row 18
$ws.Range("B18").Value = 0.4759141336213588
$ws.Range("C18").Value = 0.09517167971577578
$ws.Range("D18").Value = 0.03252563509625617
$ws.Range("E18").Value = 0.5872616221731448
$ws.Range("F18").Value = 0.6370352279284219
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.4589529743305754
$ws.Range("K18").Value = 0.5310857162884872
$ws.Range("O18").Value = 2.107484800232072
# Row 19
$ws.Range("B19").Value = 0.4720289556434807
$ws.Range("C19").Value = 0.09442896736862849
$ws.Range("D19").Value = 0.03234025375314786
$ws.Range("E19").Value = 0.5822648774036594
$ws.Range("F19").Value = 0.6367175830095562
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.4592949665982253
$ws.Range("K19").Value = 0.5266482749546242
$ws.Range("O19").Value = 2.107701283633133
# Row 20
$ws.Range("B20").Value = 0.4895095450594908
$ws.Range("C20").Value = 0.09776950116582839
$ws.Range("D20").Value = 0.03317379933169917
$ws.Range("E20").Value = 0.6047602078499636
$ws.Range("F20").Value = 0.6381790631377982
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.4577778432668396
$ws.Range("K20").Value = 0.5466118350793749
$ws.Range("O20").Value = 2.106832459512958
# Row 21
$ws.Range("B21").Value = 0.5481845301332555
$ws.Range("C21").Value = 0.1089623191430462
$ws.Range("D21").Value = 0.03596213610445176
$ws.Range("E21").Value = 0.6805135720824467
$ws.Range("F21").Value = 0.6436512441938333
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.4530641585027375
$ws.Range("K21").Value = 0.6135894279108811
$ws.Range("O21").Value = 2.105761430600808
# Row 22
$ws.Range("B22").Value = 0.5864863691370772
$ws.Range("C22").Value = 0.1162541706305831
$ws.Range("D22").Value = 0.03777513364546792
$ws.Range("E22").Value = 0.730159875476275
$ws.Range("F22").Value = 0.6476476368204089
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.4502703254260361
$ws.Range("K22").Value = 0.6572878066463659
$ws.Range("O22").Value = 2.106443748056137
# Row 23
$ws.Range("B23").Value = 0.566048318140389
$ws.Range("C23").Value = 0.1123645208459152
$ws.Range("D23").Value = 0.03680837046005792
$ws.Range("E23").Value = 0.7036495962146745
$ws.Range("F23").Value = 0.6454761739146733
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.4517351437527388
$ws.Range("K23").Value = 0.6339722381490276
$ws.Range("O23").Value = 2.1059527883192
# Row 24
$ws.Range("B24").Value = 0.4885498031998736
$ws.Range("C24").Value = 0.09758617066088959
$ws.Range("D24").Value = 0.03312807060398626
$ws.Range("E24").Value = 0.6035242465563755
$ws.Range("F24").Value = 0.6380966953545979
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.4578597154419235
$ws.Range("K24").Value = 0.5455158901696961
$ws.Range("O24").Value = 2.106873230967182
# Row 25
$ws.Range("B25").Value = 0.4047774681259853
$ws.Range("C25").Value = 0.08154691188487107
$ws.Range("D25").Value = 0.02911946681508937
$ws.Range("E25").Value = 0.4960541509792193
$ws.Range("F25").Value = 0.6319315983020317
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.4656917171671076
$ws.Range("K25").Value = 0.4497962708934722
$ws.Range("O25").Value = 2.113769074292208
